$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a Price-column value while keeping it plain text (some
# updated prices are numeric-looking, e.g. "344.00"/"60.30"; a bare
# Range.Value assignment would let Excel coerce them to numbers and drop
# the significant trailing zeros that the source sheet stores as text).
# Format-as-text, assign, then clear the temporary format so no extra
# cell style is left behind versus the original (unstyled) cells.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "60.256.14"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "2.613.76"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws.Range("D5") "522.73"
$ws.Range("E5").Value = "  +1.23%  "
Set-TextValue $ws.Range("D6") "148.91"
$ws.Range("E6").Value = "  -3.85%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -4.82%  "
$ws.Range("D9").Value = "2.618.33"
$ws.Range("E9").Value = "  +0.07%  "
Set-TextValue $ws.Range("D10") "6.32"
$ws.Range("E10").Value = "  -4.62%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "3.070.83"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "60.309.36"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("E16").Value = "  -2.13%  "
Set-TextValue $ws.Range("D17") "0.0000138"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "2.608.91"
$ws.Range("E18").Value = "  -0.17%  "
Set-TextValue $ws.Range("D19") "4.64"
$ws.Range("E19").Value = "  -2.83%  "
Set-TextValue $ws.Range("D20") "344.00"
$ws.Range("E20").Value = "  -3.61%  "
Set-TextValue $ws.Range("D21") "10.42"
$ws.Range("E21").Value = "  -1.92%  "
Set-TextValue $ws.Range("D22") "6.11"
Set-TextValue $ws.Range("D23") "0.994"
$ws.Range("E23").Value = "  -0.47%  "
Set-TextValue $ws.Range("D24") "60.30"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D27") "0.163"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").Value = "0.0₃0809"
$ws.Range("E28").Value = "  -4.09%  "
Set-TextValue $ws.Range("D29") "7.08"
$ws.Range("E29").Value = "  -3.88%  "
$ws.Range("E30").Value = "  +0.04%  "
Set-TextValue $ws.Range("D31") "6.02"
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -2.80%  "
Set-TextValue $ws.Range("D34") "149.74"
$ws.Range("E34").Value = "  +0.05%  "
Set-TextValue $ws.Range("D35") "3.98"
$ws.Range("E35").Value = "  -2.29%  "
Set-TextValue $ws.Range("D36") "0.926"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("E37").Value = "  -5.06%  "
Set-TextValue $ws.Range("D38") "0.865"
$ws.Range("E38").Value = "  +2.37%  "
Set-TextValue $ws.Range("D39") "36.51"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("E41").Value = "  -3.96%  "
Set-TextValue $ws.Range("D42") "288.07"
$ws.Range("E42").Value = "  -0.14%  "
Set-TextValue $ws.Range("D43") "0.626"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -1.06%  "
Set-TextValue $ws.Range("D45") "0.998"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("E49").Value = "  +0.98%  "
Set-TextValue $ws.Range("D50") "4.69"
$ws.Range("E50").Value = "  -5.56%  "
$ws.Range("D51").Value = "1.953.70"
$ws.Range("E51").Value = "  -1.45%  "
